$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: turn the old placeholder row into the "6.2 / Listado y detalle
#     de Componentes" entry (mirrors the "6.1" entry already in row 7). ---

# A8 & B8 contain text that LOOKS numeric-ish ("6.2") or must stay literal
# text; going through a formula + paste-values round trip keeps them as
# genuine shared-string text (t="s") without Excel's auto-number coercion
# and without minting a new cell style.
$scratch = $ws.Range("Z1")

$scratch.Formula = '="6.2"'
$scratch.Copy()
$ws.Range("A8").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B8").Value = "Listado y detalle de Componentes"
$ws.Range("C8").Value = "Jon"
$ws.Range("D8").Value = "5/4/2025"
$ws.Range("E8").Value = "5/4/2025"
$ws.Range("F8").Value = "✅ Hecho"

$ws.Rows.Item(8).RowHeight = 30

# --- Rows 9-15: the old sequential index numbers in column A are cleared
#     (the cell stays present, just empty) now that row 8 carries its own
#     "6.2" label instead of a running count. ---
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()

# --- Selection moves to G19 ---
$ws.Range("G19").Select()
